# Update the time log: finished 2 questions from 3.5 on top of the
# existing "3 problems from 3.4" entry, adding 0.5 hours of work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 67 corresponds to the JS101 entry dated 2021-10-18 that originally
# only noted "3 problems from 3.4" and logged 0.5 hours.
$ws.Range("D67").Value = "3 problems from 3.4, 2 questions from 3.5"
$ws.Range("C67").Value = 1

# The weekly total (D72 = SUM(C65:C71)) and grand total (C73 = SUBTOTAL)
# are formulas and will recalculate automatically.

# Update the active selection to reflect where the author left off.
$ws.Range("C68").Select()
